# Make IXLCell.GetFormattedString() more compliant with Excel: custom
# formats only apply to numbers, not other types (text, bool, error).
# For a boolean cell, GetFormattedString() should render the same way
# Excel does ("TRUE"/"FALSE"), not the .NET ToString() casing ("True").
#
# In the "Cell Values" sample sheet, cell G4 demonstrates
# GetFormattedString() for the boolean row (row 4) and previously held
# the text "True". Update it to "TRUE" to match Excel's convention. The
# leading apostrophe forces the value to be stored as text (not
# re-interpreted as a boolean literal) so the cell keeps its existing
# text/shared-string type and number-format style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cell Values")
$ws.Range("G4").Value = "'TRUE"
